$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the two title strings between row 2 and row 3 -----------------
$ws.Range("A2").Value = "United States presidential election of 1960"
$ws.Range("A3").Value = "U. S. Electoral College"

# --- Swap the two uri strings shown in row 2 and row 3 -------------------
$ws.Range("E2").Value = "https://www.britannica.com/event/United-States-presidential-election-of-1960"
$ws.Range("E3").Value = "https://www.archives.gov/federal-register/electoral-college/votes/1953_1957.html#1960"

# --- Re-point the hyperlinks ----------------------------------------------
# The engine's Hyperlinks.Delete() always clears every hyperlink on the
# sheet (it is not scoped to the calling range), so wipe once and rebuild
# both links with the desired end state:
#   E2 -> archives.gov (no sub-address / location anchor)
#   E3 -> britannica    (sub-address / location = "1960")
$ws.Range("E2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E2"), "https://www.archives.gov/federal-register/electoral-college/votes/1953_1957.html", "")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.britannica.com/event/United-States-presidential-election-of-1960", "1960")

# Hyperlinks.Add re-applies its own (near-duplicate) style xf; restore the
# original named "Hyperlink" cell style so E2/E3 keep their original look.
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"
